$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "element" / "weapon_type" column headers (E1 <-> F1)
$ws.Range("E1").Value = "element"
$ws.Range("F1").Value = "weapon_type"

# Mark the first 30 characters (rows 2-31) as achieved
# (leading apostrophe keeps this a literal text "true", matching the
# "true"/"false" text labels already used throughout the sheet instead of
# Excel auto-converting it to a native boolean; resetting the style afterward
# drops the "quote prefix" marker so formatting matches the rest of the sheet)
$ws.Range("B2:B31").Value = "'true"
$ws.Range("B2:B31").Style = "Normal"

# Restore the view: scroll so row 16 is at the top, select I36
$ws.Range("I36").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
